$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 43235
$ws.Range("B32").Value = "Début d'un gros refactoring de l'entierté du code"
$ws.Range("C32").Value = 3

$ws.Range("B33").Select()
